$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Wrong Answer7" -> "Wrong Answer 7" in cell I1
$ws.Range("I1").Value = "Wrong Answer 7"

# Update selection to I1
$ws.Range("I1").Select() | Out-Null
